# Update Name of Algo
# Apply updated KNN imputation results to the data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -10.878
$ws.Range("D3").Value = -7.040000000000001
$ws.Range("B12").Value = 4.807
$ws.Range("C14").Value = -11.895
$ws.Range("C26").Value = -13.04
$ws.Range("D30").Value = -7.25
$ws.Range("C31").Value = -12.575
$ws.Range("B32").Value = 6.379
$ws.Range("C35").Value = -12.311
$ws.Range("B36").Value = 8.673999999999999
$ws.Range("C37").Value = -13.387
$ws.Range("B38").Value = 5.203
$ws.Range("D44").Value = -7.409999999999999
$ws.Range("C45").Value = -12.67
$ws.Range("B46").Value = 5.823
$ws.Range("B54").Value = 5.340000000000001
$ws.Range("B55").Value = 4.742
$ws.Range("C57").Value = -13.697
$ws.Range("D58").Value = -8.021000000000001
$ws.Range("B67").Value = 5.285
$ws.Range("B69").Value = 5.034999999999999
$ws.Range("B72").Value = 5.095
$ws.Range("D84").Value = -8.238
$ws.Range("D89").Value = -6.962000000000001
$ws.Range("B91").Value = 6.218000000000001
$ws.Range("D91").Value = -6.395
$ws.Range("D92").Value = -6.681999999999999
$ws.Range("B99").Value = 5.718
$ws.Range("C100").Value = -12.807
$ws.Range("C102").Value = -13.583
$ws.Range("D102").Value = -7.723999999999999
